$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: "Road" -> "Neighbourhood " + "Road 1" (two runs, same
#    bold Calibri 24pt formatting, matching the XML diff exactly).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$rAll = $p1.Range
if ($rAll.Text -eq "Road`r") {
    $start = $rAll.Start

    # Replace "Road" (4 chars) with the combined new text.
    $rWord = $d.Range($start, $start + 4)
    $rWord.Text = "Neighbourhood Road 1"

    # Force Word to split this into two runs at the boundary between
    # "Neighbourhood " and "Road 1" by toggling Bold off/on for the second
    # portion only (it was already bold, so this is a formatting no-op that
    # still causes the run to be re-emitted as its own <w:r>).
    $splitStart = $start + "Neighbourhood ".Length
    $splitEnd = $start + "Neighbourhood Road 1".Length
    $rTail = $d.Range($splitStart, $splitEnd)
    $rTail.Font.Bold = $false
    $rTail.Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 2 & 3. Remove stray empty trailing runs (no text, just a bare
#    <w:rPr><w:rtl w:val="0"/></w:rPr> run) left over at the end of a
#    couple of paragraphs. Such runs carry no characters, so they can't be
#    reached through Range character offsets (Range.Text never shows them).
#    Instead, pull the paragraph's own OOXML via Range.XML(), strip the
#    empty-run markup with a plain string replace, and write it back with
#    Range.InsertXML so the rest of the paragraph (text, formatting, ids)
#    round-trips untouched.
# ---------------------------------------------------------------------------
$emptyRunXml = '<w:r w:rsidDel="00000000" w:rsidRPr="00000000"><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    $xml = $r.XML()
    if ($xml.Contains($emptyRunXml)) {
        $newXml = $xml.Replace($emptyRunXml, "</w:p>")
        $r.InsertXML($newXml)
    }
}

Write-Output "edits applied"
